$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 280 (shifts existing rows 280:342 down to 281:343)
$ws.Rows("280:280").Insert()

# Populate the new row 280 with the new weekly price record.
# Columns A,B,C,E,F,G,H,I,J,K,T are identical to the surrounding rows for this
# market/product block, so copy them from the row right below (now row 281).
$ws.Cells.Item(280, 1).Value2 = $ws.Cells.Item(281, 1).Value2      # A Mercado ID
$ws.Cells.Item(280, 2).Value2 = $ws.Cells.Item(281, 2).Value2      # B Mercado
$ws.Cells.Item(280, 3).Value2 = $ws.Cells.Item(281, 3).Value2      # C Region
$ws.Cells.Item(280, 4).Value2 = 44932                              # D Fecha
$ws.Cells.Item(280, 5).Value2 = $ws.Cells.Item(281, 5).Value2      # E Codreg
$ws.Cells.Item(280, 6).Value2 = $ws.Cells.Item(281, 6).Value2      # F Tipo
$ws.Cells.Item(280, 7).Value2 = $ws.Cells.Item(281, 7).Value2      # G Producto ID
$ws.Cells.Item(280, 8).Value2 = $ws.Cells.Item(281, 8).Value2      # H Producto
$ws.Cells.Item(280, 9).Value2 = $ws.Cells.Item(281, 9).Value2      # I Categoria ID
$ws.Cells.Item(280, 10).Value2 = $ws.Cells.Item(281, 10).Value2    # J Categoria
$ws.Cells.Item(280, 11).Value2 = $ws.Cells.Item(281, 11).Value2    # K Variedad
$ws.Cells.Item(280, 12).Value2 = "Primera"                         # L Calidad
$ws.Cells.Item(280, 13).Value2 = 800                                # M Volumen
$ws.Cells.Item(280, 14).Value2 = 8000                               # N Precio minimo
$ws.Cells.Item(280, 15).Value2 = 9000                               # O Precio maximo
$ws.Cells.Item(280, 16).Value2 = 8500                               # P Precio promedio ponderado
$ws.Cells.Item(280, 17).Value2 = "$/bandeja 7 kilos"                # Q Unidad de comercializacion
$ws.Cells.Item(280, 18).Value2 = "Región de La Araucanía"          # R Origen
$ws.Cells.Item(280, 19).Value2 = 1214                               # S Precio $/Kg
$ws.Cells.Item(280, 20).Value2 = $ws.Cells.Item(281, 20).Value2    # T Kg / unidad
